$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked" and "is_enabled" template header columns
# (C1/D1 held those labels; deleting them with a left shift pulls the
# remaining "order_by"/"rem" columns from E1/F1 back into C1/D1).
$ws.Range("C1:D1").Delete(-4159)
